$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("ECs", "Il17c", "Il17re", "ECs", 2, 0.6666666666666666, 1.198357333333333, 3.595072, 0.5044631819412814, 0.5044631819412814, 2, 0.6666666666666666, 0.01191533333333333, 0.035746, 0.03575915937064839, 0.03575915937064839, 0.01427882707911111, 0.128509443712, 0.01803917931966268, 0.01803917931966268),
    @("ECs", "Il17c", "Il17re", "FAPs", 2, 0.6666666666666666, 1.198357333333333, 3.595072, 0.5044631819412814, 0.5044631819412814, 3, 1, 0.154636, 0.463908, 0.4640787809914049, 0.4640787809914049, 0.1853091845973334, 1.667782661376, 0.2341106585303551, 0.2341106585303551),
    @("ECs", "Il17c", "Il17re", "sCs", 2, 0.6666666666666666, 1.198357333333333, 3.595072, 0.5044631819412814, 0.5044631819412814, 3, 1, 0.1666593333333334, 0.499978, 0.5001620596379468, 0.5001620596379468, 0.1997174342684445, 1.797456908416, 0.2523133440912636, 0.2523133440912636),
    @("FAPs", "Il17c", "Il17re", "ECs", 2, 0.6666666666666666, 0.7359279999999999, 2.207784, 0.3097978960307471, 0.309797896030747, 2, 0.6666666666666666, 0.01191533333333333, 0.035746, 0.03575915937064839, 0.03575915937064839, 0.008768827429333333, 0.07891944686399999, 0.01107811233685505, 0.01107811233685504),
    @("FAPs", "Il17c", "Il17re", "FAPs", 2, 0.6666666666666666, 0.7359279999999999, 2.207784, 0.3097978960307471, 0.309797896030747, 3, 1, 0.154636, 0.463908, 0.4640787809914049, 0.4640787809914049, 0.113800962208, 1.024208659872, 0.1437706299436511, 0.1437706299436511),
    @("FAPs", "Il17c", "Il17re", "sCs", 2, 0.6666666666666666, 0.7359279999999999, 2.207784, 0.3097978960307471, 0.309797896030747, 3, 1, 0.1666593333333334, 0.499978, 0.5001620596379468, 0.5001620596379468, 0.1226492698613333, 1.103843428752, 0.154949153750241, 0.1549491537502409),
    @("sCs", "Il17c", "Il17re", "ECs", 3, 1, 0.4412246666666667, 1.323674, 0.1857389220279715, 0.1857389220279715, 2, 0.6666666666666666, 0.01191533333333333, 0.035746, 0.03575915937064839, 0.03575915937064839, 0.005257338978222222, 0.047316050804, 0.00664186771413067, 0.006641867714130669),
    @("sCs", "Il17c", "Il17re", "FAPs", 3, 1, 0.4412246666666667, 1.323674, 0.1857389220279715, 0.1857389220279715, 3, 1, 0.154636, 0.463908, 0.4640787809914049, 0.4640787809914049, 0.06822921755466667, 0.6140629579920001, 0.08619749251739862, 0.08619749251739862),
    @("sCs", "Il17c", "Il17re", "sCs", 3, 1, 0.4412246666666667, 1.323674, 0.1857389220279715, 0.1857389220279715, 3, 1, 0.1666593333333334, 0.499978, 0.5001620596379468, 0.5001620596379468, 0.07353420879688889, 0.661807879172, 0.09289956179644225, 0.09289956179644224)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
